$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Update Coin name / Link / Volume text cells (B, C, E columns)
$ws.Range('E2').Value = '  +2.11%  '
$ws.Range('E3').Value = '  +1.94%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('E5').Value = '  +1.64%  '
$ws.Range('E6').Value = '  +2.16%  '
$ws.Range('E7').Value = '  +2.21%  '
$ws.Range('E8').Value = '  +2.02%  '
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('E10').Value = '  +15.39%  '
$ws.Range('E11').Value = '  +0.52%  '
$ws.Range('E12').Value = '  +1.09%  '
$ws.Range('E13').Value = '  +6.25%  '
$ws.Range('E14').Value = '  +1.05%  '
$ws.Range('E15').Value = '  +1.99%  '
$ws.Range('E16').Value = '  +2.32%  '
$ws.Range('E17').Value = '  -0.59%  '
$ws.Range('E18').Value = '  +4.56%  '
$ws.Range('E19').Value = '  +2.84%  '
$ws.Range('E20').Value = '  +5.54%  '
$ws.Range('E21').Value = '  +0.94%  '
$ws.Range('E22').Value = '  -0.36%  '
$ws.Range('E23').Value = '  -2.14%  '
$ws.Range('E24').Value = '  +4.33%  '
$ws.Range('E25').Value = '  -1.12%  '
$ws.Range('E26').Value = '  +0.84%  '
$ws.Range('E27').Value = '  -0.81%  '
$ws.Range('E28').Value = '  -0.25%  '
$ws.Range('E29').Value = '  +2.28%  '
$ws.Range('E30').Value = '  +2.53%  '
$ws.Range('E31').Value = '  -0.59%  '
$ws.Range('E32').Value = '  -2.28%  '
$ws.Range('E33').Value = '  +2.28%  '
$ws.Range('E34').Value = '  -0.80%  '
$ws.Range('E35').Value = '  +23.65%  '
$ws.Range('E36').Value = '  +5.85%  '
$ws.Range('E37').Value = '  +2.96%  '
$ws.Range('E38').Value = '  -0.77%  '
$ws.Range('E39').Value = '  +1.34%  '
$ws.Range('E40').Value = '  +10.32%  '
$ws.Range('E42').Value = '  +4.75%  '
$ws.Range('E43').Value = '  +5.19%  '
$ws.Range('E44').Value = '  +2.02%  '
$ws.Range('E45').Value = '  +5.52%  '
$ws.Range('E46').Value = '  -1.04%  '
$ws.Range('E47').Value = '  -0.59%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('E48').Value = '  +3.52%  '
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('E49').Value = '  +4.54%  '
$ws.Range('E50').Value = '  +0.17%  '
$ws.Range('E51').Value = '  +7.41%  '

# Update Price cells (D column) as text to preserve formatting (avoid numeric coercion)
Set-TextValue 'D2' '70.759.09'
Set-TextValue 'D3' '3.564.50'
Set-TextValue 'D5' '583.92'
Set-TextValue 'D6' '189.68'
Set-TextValue 'D8' '3.556.38'
Set-TextValue 'D12' '54.59'
Set-TextValue 'D13' '0.0000320'
Set-TextValue 'D14' '9.52'
Set-TextValue 'D15' '4.134.54'
Set-TextValue 'D16' '70.752.63'
Set-TextValue 'D17' '19.17'
Set-TextValue 'D18' '12.78'
Set-TextValue 'D19' '3.590.30'
Set-TextValue 'D20' '569.66'
Set-TextValue 'D23' '17.91'
Set-TextValue 'D25' '4.91'
Set-TextValue 'D26' '94.25'
Set-TextValue 'D27' '11.16'
Set-TextValue 'D30' '32.54'
Set-TextValue 'D31' '7.20'
Set-TextValue 'D32' '12.31'
Set-TextValue 'D34' '63.97'
Set-TextValue 'D35' '3.74'
Set-TextValue 'D36' '3.26'
Set-TextValue 'D37' '0.411'
Set-TextValue 'D38' '532.62'
Set-TextValue 'D39' '38.36'
Set-TextValue 'D40' '3.634.03'
Set-TextValue 'D41' '1.00'
Set-TextValue 'D42' '0.0₃0799'
Set-TextValue 'D45' '0.0467'
Set-TextValue 'D46' '2.95'
Set-TextValue 'D48' '0.138'
Set-TextValue 'D49' '9.28'
